$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column BF ("Date") stores a text string like "5-4-2012-13" for every
# data row. Because of how NBA.com displayed the stats, the date was off
# by one day, so every "5-4-2012-13" needs to become the literal text
# "2013-05-04" (still plain text, not an Excel date serial number).
#
# Assigning a date-shaped string straight through .Value / .Value2 /
# .Formula makes Excel auto-parse it into a real date. To write the
# literal text for real (and without permanently changing the number
# format / style of the target cells), stage the text in a scratch cell
# that's explicitly formatted as Text, copy it, and paste-special
# (values only) into each matching cell - then restore the scratch cell
# to exactly how it was before we touched it.

$oldText = "5-4-2012-13"
$newText = "2013-05-04"

$scratch = $ws.Range("A1")
$scratchOrigValue = $scratch.Value2
$scratchOrigFormat = $scratch.NumberFormat

$scratch.NumberFormat = "@"
$scratch.Value = $newText
$scratch.Copy()

for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Range("BF$r")
    if ($cell.Value2 -eq $oldText) {
        $cell.PasteSpecial(-4163)
    }
}

$excel.CutCopyMode = $false

# Put the scratch cell back exactly as found.
$scratch.Clear()
if ($scratchOrigValue -ne $null) {
    $scratch.NumberFormat = $scratchOrigFormat
    $scratch.Value = $scratchOrigValue
}
